$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $replace, 2)
}

# 1. "deux partie distincte, une bibliothèque « Donnée »" -> pluralize both nouns
Replace-Text "est divisé en deux partie distincte, une bibliothèque « Donnée »" `
             "est divisé en deux partie distinctes, une bibliothèque « Données »"

# 2. Remove " ou le tri" before the final period of the first paragraph
Replace-Text "mais aussi toutes les classes d’utilitaires, comme la recherche ou le tri." `
             "mais aussi toutes les classes d’utilitaires, comme la recherche."

# 3. Rewrite the description of what the Ensemble Audio object represents
Replace-Text "cet objet servira concrètement à représenter un album de morceau, mais aussi une playlist de podcast, où une collection de station radio. A cette clé Ensemble" `
             "cet objet servira concrètement à représenter une collection contenant des morceaux, des podcasts et des stations de radio. A cette clé Ensemble"

# 3b. Re-write across the first "LinkedList" occurrence so Word merges the runs
#     and drops the (now unwanted) spell-check proofErr markers around it.
Replace-Text "d’une LinkedList d’objet Piste" `
             "d’une LinkedList d’objet Piste"

# 3c. Same for the second "LinkedList" occurrence a bit further down.
Replace-Text "dans les LinkedList de Piste" `
             "dans les LinkedList de Piste"

# 4. "dico" -> "dictionnaire"
Replace-Text "ne seule clé du dico est sélectionnée, " `
             "ne seule clé du dictionnaire est sélectionnée, "

# 5. "infos" -> "informations", "appli" -> "application", and drop " dans le xaml"
Replace-Text " toutes les infos relatives aux paramètres de l’appli et du profil , ou encore un Manager lecteur, qui contrôle le Lecteur de Musique dans le xaml, ainsi qu’une liste de lecture associée à un album." `
             " toutes les informations relatives aux paramètres de l’application et du profil , ou encore un Manager lecteur, qui contrôle le Lecteur de Musique, ainsi qu’une liste de lecture associée à un album."

# 6. Rewrite the last paragraph about search/sort utility, singularize it and append new sentences
Replace-Text "Enfin, un utilitaire de recherche et de tri son présents, utilisés à la fois par le Manager et le code-behind, ils présentent des méthodes assez polyvalentes, qui peuvent être utilisées dans des contextes variés." `
             "Enfin, un utilitaire de recherche est présent, utilisé à la fois par le Manager et le code-behind, il présente des méthodes assez polyvalentes, qui peuvent être utilisées dans des contextes variés. On peut par exemple effectuer une recherche par genre qui n’affichera que les ensembles audio correspondant à ce genre, ou encore une recherche par mot-clé qui affichera les ensembles dont le nom correspond, ou dont une des pistes a son titre ou son artiste qui correspond à la recherche."
